# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1370
$ws1.Range("F5").Value  = 107
$ws1.Range("F7").Value  = 11720
$ws1.Range("F8").Value  = 4394
$ws1.Range("F10").Value = 41
$ws1.Range("F11").Value = 26
$ws1.Range("F13").Value = 2550
$ws1.Range("F14").Value = 1098
$ws1.Range("F16").Value = 42
$ws1.Range("F17").Value = 5109
$ws1.Range("F19").Value = 185
$ws1.Range("F21").Value = 11348
$ws1.Range("F22").Value = 11283
$ws1.Range("F24").Value = 46
$ws1.Range("F26").Value = 12
$ws1.Range("F27").Value = 48

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1370
$ws4.Range("F5").Value  = 107
$ws4.Range("F7").Value  = 11720
$ws4.Range("F8").Value  = 4395
$ws4.Range("F10").Value = 41
$ws4.Range("F11").Value = 26
$ws4.Range("F13").Value = 2550
$ws4.Range("F15").Value = 1098
$ws4.Range("F17").Value = 42
$ws4.Range("F18").Value = 5109
$ws4.Range("F20").Value = 185
$ws4.Range("F22").Value = 11348
$ws4.Range("F23").Value = 11283
$ws4.Range("F25").Value = 46
$ws4.Range("F27").Value = 12
$ws4.Range("F28").Value = 48

$wb.Save()
